# Scheduled-runner update: refresh currentAveragePrice / leve-profit figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR profit sheets.
# Generated from the upstream market-data diff - one block per sheet, one
# comment + cell-writes per affected leve row (H..N = currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road | Potion
$ws.Cells.Item(17, 8).Value = 365
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 365
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 1095
$ws.Cells.Item(17, 14).Value = -1431

# Row 113: Amaro Kart | Starch Glue
$ws.Cells.Item(113, 8).Value = 16733.1
$ws.Cells.Item(113, 9).Value = 22018.715
$ws.Cells.Item(113, 10).Value = 4400
$ws.Cells.Item(113, 11).Value = 22018.715
$ws.Cells.Item(113, 12).Value = 4400
$ws.Cells.Item(113, 13).Value = -18764.715
$ws.Cells.Item(113, 14).Value = -10908

# Row 125: Body over Mind | Grade 5 Dexterity Alkahest
$ws.Cells.Item(125, 8).Value = 2605
$ws.Cells.Item(125, 9).Value = 2600
$ws.Cells.Item(125, 10).Value = 2608.75
$ws.Cells.Item(125, 11).Value = 23400
$ws.Cells.Item(125, 12).Value = 23478.75
$ws.Cells.Item(125, 13).Value = -20940
$ws.Cells.Item(125, 14).Value = -28398.75


$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth | Bronze Rivets
$ws.Cells.Item(5, 8).Value = 51.714287
$ws.Cells.Item(5, 9).Value = 44
$ws.Cells.Item(5, 10).Value = 71
$ws.Cells.Item(5, 11).Value = 44
$ws.Cells.Item(5, 12).Value = 71
$ws.Cells.Item(5, 13).Value = 68
$ws.Cells.Item(5, 14).Value = -295

# Row 32: Ingot We Trust | Steel Ingot
$ws.Cells.Item(32, 8).Value = 4353.05
$ws.Cells.Item(32, 9).Value = 3886.2354
$ws.Cells.Item(32, 10).Value = 6998.3335
$ws.Cells.Item(32, 11).Value = 3886.2354
$ws.Cells.Item(32, 12).Value = 6998.3335
$ws.Cells.Item(32, 13).Value = -3599.2354

# Row 33: A Leg to Stand On | Heavy Iron Flanchard
$ws.Cells.Item(33, 8).Value = 506666.66
$ws.Cells.Item(33, 9).Value = 506666.66
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 506666.66
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = -506337.66

# Row 36: Hot for Teacher | Heavy Iron Armor
$ws.Cells.Item(36, 8).Value = 2023
$ws.Cells.Item(36, 9).Value = 2023
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 2023
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -1677

# Row 39: Aurochs Star | Bull Hoplon
$ws.Cells.Item(39, 8).Value = 4499
$ws.Cells.Item(39, 9).Value = 4499
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 4499
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).Value = -3979

# Row 41: Skillet Scandal | White Skillet
$ws.Cells.Item(41, 8).Value = 3150
$ws.Cells.Item(41, 9).Value = 1780
$ws.Cells.Item(41, 10).Value = 10000
$ws.Cells.Item(41, 11).Value = 1780
$ws.Cells.Item(41, 12).Value = 10000
$ws.Cells.Item(41, 13).Value = -1366
$ws.Cells.Item(41, 14).Value = -10828

# Row 43: They've Got Legs | Steel Sabatons
$ws.Cells.Item(43, 8).Value = 21799.8
$ws.Cells.Item(43, 9).Value = 20000
$ws.Cells.Item(43, 10).Value = 22249.75
$ws.Cells.Item(43, 11).Value = 20000
$ws.Cells.Item(43, 12).Value = 22249.75
$ws.Cells.Item(43, 13).Value = -19687
$ws.Cells.Item(43, 14).Value = -22875.75

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Cells.Item(61, 8).Value = 2207.5
$ws.Cells.Item(61, 9).Value = 3490
$ws.Cells.Item(61, 10).Value = 925
$ws.Cells.Item(61, 11).Value = 3490
$ws.Cells.Item(61, 12).Value = 925
$ws.Cells.Item(61, 13).Value = -3278
$ws.Cells.Item(61, 14).Value = -1349

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Cells.Item(136, 8).Value = 2207.5
$ws.Cells.Item(136, 9).Value = 3490
$ws.Cells.Item(136, 10).Value = 925
$ws.Cells.Item(136, 11).Value = 10470
$ws.Cells.Item(136, 12).Value = 2775
$ws.Cells.Item(136, 13).Value = -7920
$ws.Cells.Item(136, 14).Value = -7875

# Row 141: Essays on Equipment | Ra'Kaznar Greaves of Maiming
$ws.Cells.Item(141, 8).Value = 54000
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 54000
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 54000
$ws.Cells.Item(141, 14).Value = -64360


$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences | Bronze Rivets
$ws.Cells.Item(4, 8).Value = 51.714287
$ws.Cells.Item(4, 9).Value = 44
$ws.Cells.Item(4, 10).Value = 71
$ws.Cells.Item(4, 11).Value = 44
$ws.Cells.Item(4, 12).Value = 71
$ws.Cells.Item(4, 13).Value = 71
$ws.Cells.Item(4, 14).Value = -301

# Row 22: Riveting Run | Iron Rivets
$ws.Cells.Item(22, 8).Value = 499
$ws.Cells.Item(22, 9).Value = 499
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 499
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -326
$ws.Cells.Item(22, 14).ClearContents()


$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof | Ash Lumber
$ws.Cells.Item(16, 8).Value = 975.2
$ws.Cells.Item(16, 9).Value = 975.2
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 975.2
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -688.2

# Row 105: Zelkova, My Love | Zelkova Lumber
$ws.Cells.Item(105, 8).Value = 5344.3335
$ws.Cells.Item(105, 9).Value = 3020
$ws.Cells.Item(105, 10).Value = 8249.75
$ws.Cells.Item(105, 11).Value = 3020
$ws.Cells.Item(105, 12).Value = 8249.75
$ws.Cells.Item(105, 13).Value = -1273

# Row 113: Patient Patients | White Ash Lumber
$ws.Cells.Item(113, 8).Value = 975.2
$ws.Cells.Item(113, 9).Value = 975.2
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 975.2
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 1194.8
$ws.Cells.Item(113, 14).ClearContents()


$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap | Maple Syrup
$ws.Cells.Item(5, 8).Value = 1237.3
$ws.Cells.Item(5, 9).Value = 1451
$ws.Cells.Item(5, 10).Value = 1145.7142
$ws.Cells.Item(5, 11).Value = 4353
$ws.Cells.Item(5, 12).Value = 3437.1426
$ws.Cells.Item(5, 13).Value = -4241
$ws.Cells.Item(5, 14).Value = -3661.1426

# Row 38: Pretty as a Picture | Dark Vinegar
$ws.Cells.Item(38, 8).Value = 92.25
$ws.Cells.Item(38, 9).Value = 23
$ws.Cells.Item(38, 10).Value = 300
$ws.Cells.Item(38, 11).Value = 69
$ws.Cells.Item(38, 12).Value = 900
$ws.Cells.Item(38, 13).Value = 278
$ws.Cells.Item(38, 14).Value = -1594

# Row 135: Not-so-secret Ingredient | Royal Maple Syrup
$ws.Cells.Item(135, 8).Value = 1237.3
$ws.Cells.Item(135, 9).Value = 1451
$ws.Cells.Item(135, 10).Value = 1145.7142
$ws.Cells.Item(135, 11).Value = 13059
$ws.Cells.Item(135, 12).Value = 10311.4278
$ws.Cells.Item(135, 13).Value = -10524
$ws.Cells.Item(135, 14).Value = -15381.4278


$ws = $wb.Worksheets.Item("GSM")
# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Cells.Item(126, 8).Value = 6777.4287
$ws.Cells.Item(126, 9).Value = 5873.25
$ws.Cells.Item(126, 10).Value = 7983
$ws.Cells.Item(126, 11).Value = 17619.75
$ws.Cells.Item(126, 12).Value = 23949
$ws.Cells.Item(126, 13).Value = -15149.75
$ws.Cells.Item(126, 14).Value = -28889

# Row 132: On Board for Lar | Lar Ingot
$ws.Cells.Item(132, 8).Value = 3943.0881
$ws.Cells.Item(132, 9).Value = 3731.6667
$ws.Cells.Item(132, 10).Value = 4758.5713
$ws.Cells.Item(132, 11).Value = 11195.0001
$ws.Cells.Item(132, 12).Value = 14275.7139
$ws.Cells.Item(132, 13).Value = -8665.000100000001

# Row 136: Shiny and Good | Pink Beryl
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 14).ClearContents()


$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Cells.Item(7, 8).Value = 3250
$ws.Cells.Item(7, 9).Value = 3250
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 3250
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = -3138
$ws.Cells.Item(7, 14).ClearContents()

# Row 40: Best Served Toad | Toad Leather
$ws.Cells.Item(40, 8).Value = 4248
$ws.Cells.Item(40, 9).Value = 4248
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 4248
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -4112

# Row 55: It's Not a Job, It's a Calling | Peiste Leather
$ws.Cells.Item(55, 8).Value = 883.5
$ws.Cells.Item(55, 9).Value = 870.2
$ws.Cells.Item(55, 10).Value = 950
$ws.Cells.Item(55, 11).Value = 870.2
$ws.Cells.Item(55, 12).Value = 950
$ws.Cells.Item(55, 13).Value = -697.2
$ws.Cells.Item(55, 14).Value = -1296

# Row 68: You Could Say It's a Moving Target | Wyvern Leather
$ws.Cells.Item(68, 8).Value = 2439.2
$ws.Cells.Item(68, 9).Value = 2449
$ws.Cells.Item(68, 10).Value = 2400
$ws.Cells.Item(68, 11).Value = 2449
$ws.Cells.Item(68, 12).Value = 2400
$ws.Cells.Item(68, 13).Value = -1700

# Row 71: They Call It Bloody Mary (L) | Wyvern Leather
$ws.Cells.Item(71, 8).Value = 2439.2
$ws.Cells.Item(71, 9).Value = 2449
$ws.Cells.Item(71, 10).Value = 2400
$ws.Cells.Item(71, 11).Value = 12245
$ws.Cells.Item(71, 12).Value = 12000
$ws.Cells.Item(71, 13).Value = -8501

# Row 126: Battered Books | Saiga Leather
$ws.Cells.Item(126, 8).Value = 3250
$ws.Cells.Item(126, 9).Value = 3250
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 9750
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -7280
$ws.Cells.Item(126, 14).ClearContents()

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Cells.Item(136, 8).Value = 3995.5
$ws.Cells.Item(136, 9).Value = 4999
$ws.Cells.Item(136, 10).Value = 2992
$ws.Cells.Item(136, 11).Value = 14997
$ws.Cells.Item(136, 12).Value = 8976
$ws.Cells.Item(136, 13).Value = -12447


$ws = $wb.Worksheets.Item("WVR")
# Row 34: He's Got Legs | Velveteen Sarouel
$ws.Cells.Item(34, 8).Value = 1000
$ws.Cells.Item(34, 9).Value = 1000
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 1000
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = -797
$ws.Cells.Item(34, 14).ClearContents()

# Row 126: A Polished Purchase | Snow Linen
$ws.Cells.Item(126, 8).Value = 1007.6429
$ws.Cells.Item(126, 9).Value = 1131.2
$ws.Cells.Item(126, 10).Value = 698.75
$ws.Cells.Item(126, 11).Value = 3393.6
$ws.Cells.Item(126, 12).Value = 2096.25
$ws.Cells.Item(126, 13).Value = -923.6000000000004
$ws.Cells.Item(126, 14).Value = -7036.25

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Cells.Item(136, 8).Value = 1718.5
$ws.Cells.Item(136, 9).Value = 1718.5
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 5155.5
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -2605.5
